# C1--C2-and-C3-PowerPoint.pptx edit
#
# 1) Slide 16's summary table is re-styled to a different built-in table
#    style (swap the table style GUID applied to the table).
# 2) The deck's theme palette is switched from the "Integral" accent
#    colours over to the standard Office theme colours.

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 16 -------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{E0D07416-5522-4812-B213-E8DA77674827}")
}

# --- 2) Swap the theme colour scheme ("Integral" -> "Office") --------------
function Convert-HexToOleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# msoThemeColorDark1..msoThemeColorFollowedHyperlink (index 1-12), Office theme values
$officeThemeColors = @(
    "000000", # Dark 1
    "FFFFFF", # Light 1
    "44546A", # Dark 2
    "E7E6E6", # Light 2
    "5B9BD5", # Accent 1
    "ED7D31", # Accent 2
    "A5A5A5", # Accent 3
    "FFC000", # Accent 4
    "4472C4", # Accent 5
    "70AD47", # Accent 6
    "0563C1", # Hyperlink
    "954F72"  # Followed Hyperlink
)

$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColorScheme.Item($i).RGB = Convert-HexToOleColor($officeThemeColors[$i - 1])
}
